# Trade #129 closed at 2026-02-18 00:40:48 - unknown UNKNOWN +0.000%
#
# This script reproduces, via the Excel COM object model, the changes that
# were applied to live_trading_results.xlsx:
#   - Summary KPIs refreshed (capital / P&L / trade counts / win rate)
#   - Strategy Status row for HighProbConvergence refreshed
#   - Trade #157 (HighProbConvergence) moved from OPEN to CLOSED (early_exit)
#     on both the "All Trades" ledger and its strategy-specific sheet
#   - Two brand new OPEN trades appended (#186 momentum, #187 MarketMaking)
#     to "All Trades" and to their respective strategy sheets

$wb = $excel.ActiveWorkbook

function Set-DateLikeTextCell {
    # Writes a literal text value that LOOKS like a date/time (e.g.
    # "2026-02-18" or "00:40:42") into a cell without letting Excel's
    # automatic cell-content parsing reinterpret it as a real date/time
    # serial number, then drops the temporary text-number-format so the
    # cell is left with the default (unstyled) look, matching how the
    # ledger stores these columns as plain strings.
    param(
        $Sheet,
        [int]$Row,
        [int]$Col,
        [string]$Text
    )
    $cell = $Sheet.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $Text
    $cell.ClearFormats()
}

# ---------------------------------------------------------------------
# Summary
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.09   # Current Capital
$summary.Range("B4").Value = 0.2       # Total P&L $
$summary.Range("B5").Value = 0.03      # Total P&L %
$summary.Range("B6").Value = 157       # Total Trades
$summary.Range("B7").Value = 71        # Winning Trades
$summary.Range("B9").Value = 45.22     # Win Rate %

# ---------------------------------------------------------------------
# Strategy Status - HighProbConvergence row (row 3)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C3").Value = 100.38     # Capital
$status.Range("D3").Value = 20         # Trades
$status.Range("E3").Value = 0.39       # P&L $
$status.Range("F3").Value = 0.38       # P&L %
$status.Range("G3").Value = 65         # Win Rate %

# ---------------------------------------------------------------------
# All Trades - Trade #157 (HighProbConvergence) closes out
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Cells.Item(158, 7).Value = 0.78          # G158 Exit Price
$allTrades.Cells.Item(158, 8).Value = "CLOSED"      # H158 Status
$allTrades.Cells.Item(158, 9).Value = 1.2987        # I158 P&L %
$allTrades.Cells.Item(158, 10).Value = 0.01         # J158 P&L $
$allTrades.Cells.Item(158, 11).Value = 100.38       # K158 Capital After
$allTrades.Cells.Item(158, 12).Value = "early_exit" # L158 Exit Reason
$allTrades.Cells.Item(158, 13).Value = 0.21         # M158 Duration (min)

# ---------------------------------------------------------------------
# All Trades - new row 187: Trade #186 (momentum, OPEN)
# ---------------------------------------------------------------------
$allTrades.Cells.Item(187, 1).Value = 186                      # A Trade #
Set-DateLikeTextCell $allTrades 187 2 "2026-02-18"              # B Date
Set-DateLikeTextCell $allTrades 187 3 "00:40:42"                # C Time
$allTrades.Cells.Item(187, 4).Value = "momentum"                # D Strategy
$allTrades.Cells.Item(187, 5).Value = "DOWN"                    # E Side
$allTrades.Cells.Item(187, 6).Value = 0.77                      # F Entry Price
# G Exit Price left blank (trade still OPEN)
$allTrades.Cells.Item(187, 8).Value = "OPEN"                    # H Status
$allTrades.Cells.Item(187, 9).Value = 0                         # I P&L %
$allTrades.Cells.Item(187, 10).Value = 0                        # J P&L $
$allTrades.Cells.Item(187, 11).Value = 99.14712996249175        # K Capital After
# L Exit Reason left blank (trade still OPEN)
$allTrades.Cells.Item(187, 13).Value = 0                        # M Duration (min)
$allTrades.Cells.Item(187, 14).Value = 0                        # N Entry Slippage (bps)
$allTrades.Cells.Item(187, 15).Value = 0                        # O Exit Slippage (bps)
$allTrades.Cells.Item(187, 16).Value = 0.9                      # P Confidence
$allTrades.Cells.Item(187, 17).Value = "Downward momentum: -45.109% over 10 samples"  # Q Entry Reason

# ---------------------------------------------------------------------
# All Trades - new row 188: Trade #187 (MarketMaking, OPEN)
# ---------------------------------------------------------------------
$allTrades.Cells.Item(188, 1).Value = 187                       # A Trade #
Set-DateLikeTextCell $allTrades 188 2 "2026-02-18"               # B Date
Set-DateLikeTextCell $allTrades 188 3 "00:40:43"                 # C Time
$allTrades.Cells.Item(188, 4).Value = "MarketMaking"             # D Strategy
$allTrades.Cells.Item(188, 5).Value = "UP"                       # E Side
$allTrades.Cells.Item(188, 6).Value = 0.22                       # F Entry Price
# G Exit Price left blank (trade still OPEN)
$allTrades.Cells.Item(188, 8).Value = "OPEN"                     # H Status
$allTrades.Cells.Item(188, 9).Value = 0                          # I P&L %
$allTrades.Cells.Item(188, 10).Value = 0                         # J P&L $
$allTrades.Cells.Item(188, 11).Value = 99.28858346467945         # K Capital After
# L Exit Reason left blank (trade still OPEN)
$allTrades.Cells.Item(188, 13).Value = 0                         # M Duration (min)
$allTrades.Cells.Item(188, 14).Value = 0                         # N Entry Slippage (bps)
$allTrades.Cells.Item(188, 15).Value = 0                         # O Exit Slippage (bps)
$allTrades.Cells.Item(188, 16).Value = 0.6                       # P Confidence
$allTrades.Cells.Item(188, 17).Value = "Normal spread capture: 198 bps"  # Q Entry Reason

# ---------------------------------------------------------------------
# momentum sheet - new row 49: Trade #186 (OPEN)
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Cells.Item(49, 1).Value = 186                          # A Trade #
Set-DateLikeTextCell $momentum 49 2 "2026-02-18"                  # B Date
Set-DateLikeTextCell $momentum 49 3 "00:40:42"                    # C Time
$momentum.Cells.Item(49, 4).Value = "momentum"                    # D Strategy
$momentum.Cells.Item(49, 5).Value = "DOWN"                        # E Side
$momentum.Cells.Item(49, 6).Value = 0.77                          # F Entry Price
# G Exit Price left blank (trade still OPEN)
$momentum.Cells.Item(49, 8).Value = "OPEN"                        # H Status
$momentum.Cells.Item(49, 9).Value = 0                             # I P&L %
$momentum.Cells.Item(49, 10).Value = 0                            # J P&L $
$momentum.Cells.Item(49, 11).Value = 99.14712996249175            # K Capital After
$momentum.Cells.Item(49, 12).Value = 0                            # L Entry Slippage (bps)
$momentum.Cells.Item(49, 13).Value = 0                            # M Exit Slippage (bps)
$momentum.Cells.Item(49, 14).Value = 0.9                          # N Confidence
$momentum.Cells.Item(49, 15).Value = "Downward momentum: -45.109% over 10 samples"  # O Entry Reason
# P Exit Reason left blank (trade still OPEN)
$momentum.Cells.Item(49, 17).Value = 0                            # Q Duration (min)

# ---------------------------------------------------------------------
# HighProbConvergence sheet - Trade #157 closes out (row 21)
# ---------------------------------------------------------------------
$hpc = $wb.Worksheets.Item("HighProbConvergence")
$hpc.Cells.Item(21, 7).Value = 0.78                 # G21 Exit Price
$hpc.Cells.Item(21, 8).Value = "CLOSED"             # H21 Status
$hpc.Cells.Item(21, 9).Value = 1.2987               # I21 P&L %
$hpc.Cells.Item(21, 10).Value = 0.01                # J21 P&L $
$hpc.Cells.Item(21, 11).Value = 100.38              # K21 Capital After
$hpc.Cells.Item(21, 16).Value = "early_exit"        # P21 Exit Reason
$hpc.Cells.Item(21, 17).Value = 0.21                # Q21 Duration (min)

# ---------------------------------------------------------------------
# MarketMaking sheet - new row 79: Trade #187 (OPEN)
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Cells.Item(79, 1).Value = 187                                # A Trade #
Set-DateLikeTextCell $mm 79 2 "2026-02-18"                        # B Date
Set-DateLikeTextCell $mm 79 3 "00:40:43"                          # C Time
$mm.Cells.Item(79, 4).Value = "MarketMaking"                      # D Strategy
$mm.Cells.Item(79, 5).Value = "UP"                                # E Side
$mm.Cells.Item(79, 6).Value = 0.22                                # F Entry Price
# G Exit Price left blank (trade still OPEN)
$mm.Cells.Item(79, 8).Value = "OPEN"                              # H Status
$mm.Cells.Item(79, 9).Value = 0                                   # I P&L %
$mm.Cells.Item(79, 10).Value = 0                                  # J P&L $
$mm.Cells.Item(79, 11).Value = 99.28858346467945                  # K Capital After
$mm.Cells.Item(79, 12).Value = 0                                  # L Entry Slippage (bps)
$mm.Cells.Item(79, 13).Value = 0                                  # M Exit Slippage (bps)
$mm.Cells.Item(79, 14).Value = 0.6                                # N Confidence
$mm.Cells.Item(79, 15).Value = "Normal spread capture: 198 bps"   # O Entry Reason
# P Exit Reason left blank (trade still OPEN)
$mm.Cells.Item(79, 17).Value = 0                                  # Q Duration (min)

Write-Host "Applied live_trading_results update for trade #157 close + trades #186/#187 open."
